# Update CodeSystem-encounter-class.xlsx metadata/concepts to new release.

$wb = $excel.ActiveWorkbook
$metaSheet = $wb.Worksheets.Item("Metadata")
$conceptSheet = $wb.Worksheets.Item("Concepts")

# 1. Bump version number
$metaSheet.Range("B3").Value = "6.0.0"

# 2. Update publication date
$metaSheet.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# 3. Fill in Publisher value
$metaSheet.Range("B9").Value = "Alvearie Team"

# 4. Replace the "Contact" row with "Jurisdiction" / "United States of America"
$metaSheet.Range("A10").Value = "Jurisdiction"
$metaSheet.Range("B10").Value = "United States of America"

# 5. Remove the duplicate "Contact" row (old row 11); this shifts remaining
#    rows (Description, Purpose, Copyright, Case Sensitive, ...) up by one.
$metaSheet.Rows.Item(11).Delete()

# 6. Set "Case Sensitive" value to the literal text "true" (now at row 14
#    after the deletion). Assigning the plain string via .Value would be
#    auto-coerced to an Excel Boolean, so build it as a text formula first
#    and then paste-special it back as a static value to keep it a string.
$metaSheet.Range("B14").Formula = '="true"'
$metaSheet.Range("B14").Copy()
$metaSheet.Range("B14").PasteSpecial(-4163)

# 7. Swap the order of the OTHER / UNKNOWN concept rows on the Concepts sheet
$conceptSheet.Range("B2").Value = "OTHER"
$conceptSheet.Range("C2").Value = "other"
$conceptSheet.Range("D2").Value = "Class of the encounter does not fit other classifications"

$conceptSheet.Range("B3").Value = "UNKNOWN"
$conceptSheet.Range("C3").Value = "unknown"
$conceptSheet.Range("D3").Value = "Class of the encounter is unknown or unspecified"
